$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Control 46
$ws.Range("D2").Value = 0.0005218258899696973
$ws.Range("E2").Value = 0.0005218258899696973

# Row 3 - Control 28
$ws.Range("D3").Value = 0.8345719029840536
$ws.Range("E3").Value = 0.8345719029840536

# Row 4 - Control 13
$ws.Range("D4").Value = 0.09030762985731447
$ws.Range("E4").Value = 0.09030762985731447

# Row 5 - Control 50
$ws.Range("D5").Value = [double]"2.449656223581634E-28"
$ws.Range("E5").Value = [double]"2.449656223581634E-28"

# Row 6 - Control 51
$ws.Range("D6").Value = 0.02053818783066794
$ws.Range("E6").Value = 0.02053818783066794

# Row 7 - MDD 35
$ws.Range("D7").Value = 0.01813808071342189
$ws.Range("E7").Value = 0.9818619192865781

# Row 8 - MDD 22
$ws.Range("D8").Value = 0.9999999298597084
$ws.Range("E8").Value = [double]"7.014029157481616E-08"

# Row 9 - MDD 50
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = 0.8684868988084139
$ws.Range("E9").Value = 0.1315131011915861

# Row 10 - MDD 45
$ws.Range("D10").Value = [double]"3.23550698390391E-14"
$ws.Range("E10").Value = 0.9999999999999677

# Row 11 - MDD 28
$ws.Range("D11").Value = 0.9998534121677101
$ws.Range("E11").Value = 0.000146587832289935
$ws.Range("F11").Value = 3.712804079055786
$ws.Range("G11").Value = 0.7
